# Apply the "trafos" sheet change: insert a new "v_base_kV" column
# (with value 132) right after "V_lv_kV" (column D), shifting the
# existing V_SCH_pu..tap_max columns one place to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("trafos")

# Shift existing header row (1) and data row (2) values from columns
# E..N to F..O, working right-to-left so we never clobber a value
# before it has been copied.
$srcCols = @("N", "M", "L", "K", "J", "I", "H", "G", "F", "E")
$dstCols = @("O", "N", "M", "L", "K", "J", "I", "H", "G", "F")

for ($i = 0; $i -lt $srcCols.Length; $i++) {
    $src = $srcCols[$i]
    $dst = $dstCols[$i]
    $ws.Range($dst + "1").Value = $ws.Range($src + "1").Value()
    $ws.Range($dst + "2").Value = $ws.Range($src + "2").Value()
}

# Insert the new column's header + value.
$ws.Range("E1").Value = "v_base_kV"
$ws.Range("E2").Value = 132

# Update the selection to cover the new used range, matching the saved file.
[void]$ws.Range("A1:O2").Select()

Write-Host "trafos sheet updated"
